$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.368.28"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "3.619.70"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  -0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.78"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.213"
$ws.Range("E9").Value = "  +7.19%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.645"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.22"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "4.194.09"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "606.03"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.00"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").Value = "70.471.16"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "3.618.03"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.02"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.10"
$ws.Range("E22").Value = "  -1.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.21"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.24"
$ws.Range("E24").Value = "  +1.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.61"
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.97"
$ws.Range("E26").Value = "  -7.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.60"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.70"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.79"
$ws.Range("E29").Value = "  +1.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.69"
$ws.Range("E30").Value = "  +7.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.26"
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.116"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.28"
$ws.Range("E34").Value = "  +0.22%  "
$ws.Range("D35").Value = "0.0₃0890"
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("D36").Value = "3.930.60"
$ws.Range("E36").Value = "  +5.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "523.49"
$ws.Range("E37").Value = "  +7.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.66"
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("E42").Value = "  -2.51%  "
$ws.Range("E43").Value = "  +2.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0461"
$ws.Range("E44").Value = "  +1.57%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.51"
$ws.Range("E45").Value = "  +6.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("E51").Value = "  +0.47%  "
